{"js": "// The diff inserts two new paragraphs (both styled \"Normal\") right after the\n// first empty paragraph that follows the long \"CummeRbund is an R package...\"\n// description paragraph:\n//   - a new, empty paragraph\n//   - a new paragraph containing the text \"Reference Style\"\n// The two originally-trailing empty paragraphs end up pushed further down,\n// unchanged.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the description paragraph (\"...Cuffdiff 2...\") so the insertion\n// point is found by content rather than a hard-coded index.\nlet descriptionIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Cuffdiff 2\") !== -1) {\n    descriptionIndex = i;\n    break;\n  }\n}\n\nif (descriptionIndex === -1) {\n  throw new Error('Could not locate the \"Cuffdiff 2\" description paragraph.');\n}\n\n// The paragraph immediately after the description is the first of the\n// trailing empty paragraphs - that's where the new content goes.\nconst anchorParagraph = paragraphs.items[descriptionIndex + 1];\n\n// Insert a new empty paragraph after the anchor, then a \"Reference Style\"\n// paragraph after that new empty paragraph. Both inherit the surrounding\n// \"Normal\" paragraph style, matching the rest of the document.\nconst newEmptyParagraph = anchorParagraph.insertParagraph(\"\", Word.InsertLocation.after);\nnewEmptyParagraph.insertParagraph(\"Reference Style\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# The diff inserts two new paragraphs (both styled \"Normal\") right after the\n# first empty paragraph that follows the long \"CummeRbund is an R package...\"\n# description paragraph:\n#   - a new, empty paragraph\n#   - a new paragraph containing the text \"Reference Style\"\n# The two originally-trailing empty paragraphs end up pushed further down,\n# unchanged.\n\n$d = $word.ActiveDocument\n\n# Locate the description paragraph (\"...Cuffdiff 2...\") by content rather\n# than a hard-coded index.\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs($i).Range.Text -match \"Cuffdiff 2\") {\n        # The paragraph right after the description is the first of the\n        # trailing empty paragraphs - that's where the new content goes.\n        $anchorIndex = $i + 1\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw 'Could not locate the \"Cuffdiff 2\" description paragraph.'\n}\n\n# Insert a new empty paragraph after the anchor, then a \"Reference Style\"\n# paragraph after that new empty paragraph. Re-fetch paragraphs by (fresh)\n# index each time, since newly inserted paragraphs shift the collection.\n$d.Paragraphs($anchorIndex).Range.InsertParagraphAfter()\n$d.Paragraphs($anchorIndex + 1).Range.InsertParagraphAfter()\n$d.Paragraphs($anchorIndex + 2).Range.Text = \"Reference Style\"\n"}
